$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot" -> add column CF (05-sep) ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Copy the header formatting (bold font, borders, centered style) from CE1
# into CF1, then set its text so the new header matches the look of the
# other day-headers.
$ws1.Range("CE1").Copy($ws1.Range("CF1"))
$ws1.Range("CF1").Value = "05-sep"

# Fill the hourly price values for 05-sep (CF2:CF25)
$ws1.Range("CF2").Value = 70.86
$ws1.Range("CF3").Value = 58.98
$ws1.Range("CF4").Value = 61.49
$ws1.Range("CF5").Value = 50.12
$ws1.Range("CF6").Value = 52.65
$ws1.Range("CF7").Value = 52.53
$ws1.Range("CF8").Value = 76.79000000000001
$ws1.Range("CF9").Value = 89.98999999999999
$ws1.Range("CF10").Value = 88.22
$ws1.Range("CF11").Value = 79.09999999999999
$ws1.Range("CF12").Value = 35
$ws1.Range("CF13").Value = 10
$ws1.Range("CF14").Value = 15
$ws1.Range("CF15").Value = 18.43
$ws1.Range("CF16").Value = 6.62
$ws1.Range("CF17").Value = 17.07
$ws1.Range("CF18").Value = 18.63
$ws1.Range("CF19").Value = 30.4
$ws1.Range("CF20").Value = 45.23
$ws1.Range("CF21").Value = 89.25
$ws1.Range("CF22").Value = 103.97
$ws1.Range("CF23").Value = 107.1
$ws1.Range("CF24").Value = 95
$ws1.Range("CF25").Value = 88.38

# --- Sheet "Gaz" -> add row 81 (2025-09-03, 30.875) ---
$ws2 = $wb.Worksheets.Item("Gaz")

# Force the new date cell to be stored as text (like every other date cell
# in column A) instead of being auto-converted to a date serial number.
$ws2.Range("A81").NumberFormat = "@"
$ws2.Range("A81").Value = "2025-09-03"
# Re-apply the (unstyled) look of the cell above so we don't leave a
# left-over "text" number format on the new cell.
$ws2.Range("A80").Copy()
$ws2.Range("A81").PasteSpecial(-4122)

$ws2.Range("B81").Value = 30.875

# --- Sheet "CO2" -> add row 81 (2025-09-03, blank value) ---
$ws3 = $wb.Worksheets.Item("CO2")

$ws3.Range("A81").NumberFormat = "@"
$ws3.Range("A81").Value = "2025-09-03"
$ws3.Range("A80").Copy()
$ws3.Range("A81").PasteSpecial(-4122)

# B81 stays blank (no CO2 price recorded yet for this date), mirroring the
# other blank cells already present in this column (e.g. B52, B69:B71).
$ws3.Range("B52").Copy()
$ws3.Range("B81").PasteSpecial(-4122)
